$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.112.98"
$ws.Range("E2").Value = "'  -0.25%  "

$ws.Range("D3").Value = "'1.875.50"
$ws.Range("E3").Value = "'  -1.88%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.10%  "

$ws.Range("D5").Value = "'313.40"
$ws.Range("E5").Value = "'  -0.41%  "

$ws.Range("E6").Value = "'  +0.09%  "

$ws.Range("D7").Value = "'0.5051"
$ws.Range("E7").Value = "'  -0.40%  "

$ws.Range("D8").Value = "'0.3846"
$ws.Range("E8").Value = "'  -2.09%  "

$ws.Range("D9").Value = "'0.08630"
$ws.Range("E9").Value = "'  -7.41%  "

$ws.Range("D10").Value = "'1.115"
$ws.Range("E10").Value = "'  -2.29%  "

$ws.Range("B11").Value = "'Polkadot"
$ws.Range("C11").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D11").Value = "'6.318"
$ws.Range("E11").Value = "'  -1.26%  "

$ws.Range("B12").Value = "'Solana"
$ws.Range("C12").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'20.67"
$ws.Range("E12").Value = "'  -1.25%  "

$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.874.30"
$ws.Range("E13").Value = "'  -2.39%  "

$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.200"
$ws.Range("E14").Value = "'  -1.76%  "

$ws.Range("B15").Value = "'BinanceUSD"
$ws.Range("C15").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "'  +0.10%  "

$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001100"
$ws.Range("E16").Value = "'  -2.03%  "

$ws.Range("B17").Value = "'Litecoin"
$ws.Range("C17").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'91.00"
$ws.Range("E17").Value = "'  -1.59%  "

$ws.Range("B18").Value = "'TRON"
$ws.Range("C18").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06631"
$ws.Range("E18").Value = "'  +0.17%  "

$ws.Range("B19").Value = "'Avalanche"
$ws.Range("C19").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'18.07"
$ws.Range("E19").Value = "'  +0.43%  "

$ws.Range("B20").Value = "'Dai"
$ws.Range("C20").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "'  +0.19%  "

$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.096"
$ws.Range("E21").Value = "'  -2.08%  "

$ws.Range("B22").Value = "'WrappedBTC"
$ws.Range("C22").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "'28.141.16"
$ws.Range("E22").Value = "'  -0.35%  "

$ws.Range("B23").Value = "'Cosmos"
$ws.Range("C23").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'11.41"
$ws.Range("E23").Value = "'  -1.22%  "

$ws.Range("B24").Value = "'Toncoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.259"
$ws.Range("E24").Value = "'  -2.85%  "

$ws.Range("B25").Value = "'LidoDAOToken"
$ws.Range("C25").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.580"
$ws.Range("E25").Value = "'  -0.67%  "

$ws.Range("B26").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "'2.091.31"
$ws.Range("E26").Value = "'  -2.23%  "

$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.75"
$ws.Range("E27").Value = "'  -1.91%  "

$ws.Range("B28").Value = "'Monero"
$ws.Range("C28").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'157.15"
$ws.Range("E28").Value = "'  -0.57%  "

$ws.Range("B29").Value = "'BitcoinCash"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'126.37"
$ws.Range("E29").Value = "'  -0.64%  "

$ws.Range("B30").Value = "'Stellar"
$ws.Range("C30").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1055"
$ws.Range("E30").Value = "'  -1.82%  "

$ws.Range("B31").Value = "'ImmutableX"
$ws.Range("C31").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.062"
$ws.Range("E31").Value = "'  -4.00%  "

$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.617"
$ws.Range("E32").Value = "'  -0.83%  "

$ws.Range("B33").Value = "'HuobiToken"
$ws.Range("C33").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.592"
$ws.Range("E33").Value = "'  -0.51%  "

$ws.Range("B34").Value = "'FraxShare"
$ws.Range("C34").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D34").Value = "'9.630"
$ws.Range("E34").Value = "'  -0.98%  "

$ws.Range("B35").Value = "'VeChain"
$ws.Range("C35").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02446"
$ws.Range("E35").Value = "'  +0.40%  "

$ws.Range("B36").Value = "'Hedera"
$ws.Range("C36").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06589"
$ws.Range("E36").Value = "'  -1.48%  "

$ws.Range("B37").Value = "'Algorand"
$ws.Range("C37").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2182"
$ws.Range("E37").Value = "'  -1.27%  "

$ws.Range("B38").Value = "'ARBITRUM"
$ws.Range("C38").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.211"
$ws.Range("E38").Value = "'  -2.58%  "

$ws.Range("B39").Value = "'TrustWalletToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.242"
$ws.Range("E39").Value = "'  -3.25%  "

$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6390"
$ws.Range("E40").Value = "'  -2.06%  "

$ws.Range("B41").Value = "'Aptos"
$ws.Range("C41").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'11.51"
$ws.Range("E41").Value = "'  -0.51%  "

$ws.Range("B42").Value = "'InternetComputer(DFINITY)"
$ws.Range("C42").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "'4.901"
$ws.Range("E42").Value = "'  -2.44%  "

$ws.Range("B43").Value = "'Frax"
$ws.Range("C43").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "'  +0.06%  "

$ws.Range("B44").Value = "'Decentraland"
$ws.Range("C44").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6006"
$ws.Range("E44").Value = "'  -2.16%  "

$ws.Range("D45").Value = "'13.17"
$ws.Range("E45").Value = "'  -1.93%  "

$ws.Range("B46").Value = "'WEMIXTOKEN"
$ws.Range("C46").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.281"
$ws.Range("E46").Value = "'  -0.60%  "

$ws.Range("B47").Value = "'PancakeSwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.674"
$ws.Range("E47").Value = "'  -1.34%  "

$ws.Range("B48").Value = "'EOS"
$ws.Range("C48").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").Value = "'1.225"
$ws.Range("E48").Value = "'  +2.99%  "

$ws.Range("D49").Value = "'1.990"
$ws.Range("E49").Value = "'  -1.76%  "

$ws.Range("B50").Value = "'Quant"
$ws.Range("C50").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'121.63"
$ws.Range("E50").Value = "'  -0.89%  "

$ws.Range("B51").Value = "'Aave"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'80.55"
$ws.Range("E51").Value = "'  +2.67%  "
